# Applies the "cleanup and go-trough first notebooks" edit to the
# "Setup workshop" slide:
#   1. Moves the existing numbered-steps textbox ("Tekstvak 5") down/left a bit.
#   2. Adds a new "Instruction steps" textbox above it, whose text is a
#      hyperlink.
#
# EMU (English Metric Units) are the native unit PowerPoint stores in the
# OOXML (914400 EMU = 1 inch = 72 points). The COM object model only takes
# points, and the runtime's point->EMU conversion is done in (lossy)
# single-precision float, which occasionally truncates a whisker below the
# intended integer EMU value. Nudging by a tiny fraction of a point keeps
# the value on the correct side of the rounding boundary without being
# large enough to ever overshoot into the next EMU.
function EmuToPt($emu) {
    return ($emu / 914400.0 * 72.0) + 0.000039
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# --- 1. Reposition the existing "Tekstvak 5" steps textbox -----------------
$steps = $s.Shapes.Item(2)
$steps.Left = EmuToPt 1284401
$steps.Top  = EmuToPt 3324496

# --- 2. Add the new "Instruction steps" hyperlinked textbox ----------------
$instr = $s.Shapes.AddTextbox(
    1,
    (EmuToPt 1143000),
    (EmuToPt 1603693),
    (EmuToPt 6094428),
    (EmuToPt 400110)
)
$instr.Name = "Tekstvak 7"

$instr.TextFrame.WordWrap = $true
$instr.TextFrame.AutoSize = 1
$instr.Fill.Visible = $false

$tr = $instr.TextFrame.TextRange
$tr.Text = "Instruction steps"
$tr.Font.Size = 20
$tr.LanguageID = "nl-BE"
$tr.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/simonperneel/workshop-UWC"
